# Generate Report for Handoff
# b.md has moved from "Handed back: in sync with en-US" to "Ready for handoff":
# a new handoff round-trip (xlf files + timestamps) was produced for b.md in
# both the zh-cn and de-de locales; the Overview sheet rolls the new handoff
# date up for b.md as well.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$overviewDate = "2016-03-24 00:35:17"
$zhHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate = "2016-03-24 00:35:13"
$deHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate = "2016-03-24 00:35:17"

# ---------------------------------------------------------------------------
# Overview sheet: row 3 (b.md) status + handoff date columns
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $overviewDate

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) status / latest handoff file & datetime, plus the
# hyperlink display text on D3 must track the new handoff file name.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = $zhHandoffFile
$wsZh.Range("E3").Value = $zhHandoffDate

$zhLinks = @()
foreach ($h in $wsZh.Hyperlinks) {
    $zhLinks += , @($h.Range.Address(), $h.Address(), $h.TextToDisplay())
}
$wsZh.Range("A1").Hyperlinks.Delete()
foreach ($link in $zhLinks) {
    $addr = $link[0]
    $target = $link[1]
    $display = $link[2]
    if ($addr -eq "`$D`$3") {
        $display = $zhHandoffFile
    }
    $wsZh.Hyperlinks.Add($wsZh.Range($addr), $target, "", "", $display)
}

# ---------------------------------------------------------------------------
# de-de sheet: row 3 (b.md) status / latest handoff file & datetime, plus the
# hyperlink display text on D3 must track the new handoff file name.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = $deHandoffFile
$wsDe.Range("E3").Value = $deHandoffDate

$deLinks = @()
foreach ($h in $wsDe.Hyperlinks) {
    $deLinks += , @($h.Range.Address(), $h.Address(), $h.TextToDisplay())
}
$wsDe.Range("A1").Hyperlinks.Delete()
foreach ($link in $deLinks) {
    $addr = $link[0]
    $target = $link[1]
    $display = $link[2]
    if ($addr -eq "`$D`$3") {
        $display = $deHandoffFile
    }
    $wsDe.Hyperlinks.Add($wsDe.Range($addr), $target, "", "", $display)
}
